$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in rows 91-110 with Lab 9 log entries ---

# Dates (column A), formatted like existing date cells (copy format from A90)
$ws.Range("A91:A110").Value = 42808
$ws.Range("A90").Copy()
$ws.Range("A91:A110").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Change/Action (B) and Who-will-it-affect (E) text, entered in original authoring order
$ws.Range("B91").Value = "Create Accounting Resources group"
$ws.Range("B92").Value = "Add Accounting group to accounting resources group"
$ws.Range("B93").Value = "Create Research Resources group"
$ws.Range("B94").Value = "Add Research-Dev group to Research Resources group"
$ws.Range("B95").Value = "Create Sales Resources domain group"
$ws.Range("B96").Value = "Add Sales group to Sales Resources group"
$ws.Range("B97").Value = "Create Departments, Personnel & Finances folders"
$ws.Range("B99").Value = "Create departmental folders Accounting,Research,Sales, & Support within Departments folder"
$ws.Range("E92").Value = "Accounting"
$ws.Range("E94").Value = "Research-Dev"
$ws.Range("E96").Value = "Sales"
$ws.Range("E103").Value = "Support"
$ws.Range("B98").Value = "Share Personnel & Finances folders"
$ws.Range("B100").Value = "Share Accounting folder, Set Accounting Resources NTFS Permissions `"Full Control`""
$ws.Range("B101").Value = "Share Research folder, Set Research Resources NTFS Permissions `"Full Control`""
$ws.Range("B102").Value = "Share Sales folder, Set Sales Resources NTFS Permissions `"Full Control`""
$ws.Range("B103").Value = "Share Support folder, Set Support Resources NTFS Permissions `"Full Control`""
$ws.Range("B104").Value = "Change Managers Domain Distribution group to Global Security"
$ws.Range("E104").Value = "Managers"
$ws.Range("B105").Value = "Set Personnel NTFS File permissions to `"Full Control`" by Managers group, Removed inherited permissions"
$ws.Range("B106").Value = "Create shared folders Projects, Shared, Archives, Confidential, & Components"
$ws.Range("B107").Value = "Share Users folder with names Users & Home"
$ws.Range("B108").Value = "Enable caching of Archives folder & Optimize for performance"
$ws.Range("B109").Value = "Disable caching of Confidential shared folder"
$ws.Range("E110").Value = "Users"
$ws.Range("B110").Value = "Grant Users group Modify and Write access to Shared folder in NTFS permissions"

# Down Time? (C), How Long (D), Tech (F), Approval (G) - reuse existing values
$ws.Range("C91").Value = "No"
$ws.Range("D91").Value = "N/A"
$ws.Range("F91").Value = "Evan"
$ws.Range("G91").Value = "ES"
$ws.Range("C92").Value = "No"
$ws.Range("D92").Value = "N/A"
$ws.Range("F92").Value = "Evan"
$ws.Range("G92").Value = "ES"
$ws.Range("C93").Value = "No"
$ws.Range("D93").Value = "N/A"
$ws.Range("F93").Value = "Evan"
$ws.Range("G93").Value = "ES"
$ws.Range("C94").Value = "No"
$ws.Range("D94").Value = "N/A"
$ws.Range("F94").Value = "Evan"
$ws.Range("G94").Value = "ES"
$ws.Range("C95").Value = "No"
$ws.Range("D95").Value = "N/A"
$ws.Range("F95").Value = "Evan"
$ws.Range("C96").Value = "No"
$ws.Range("D96").Value = "N/A"
$ws.Range("F96").Value = "Evan"
$ws.Range("G96").Value = "ES"
$ws.Range("C97").Value = "No"
$ws.Range("D97").Value = "N/A"
$ws.Range("F97").Value = "Evan"
$ws.Range("G97").Value = "ES"
$ws.Range("C98").Value = "No"
$ws.Range("D98").Value = "N/A"
$ws.Range("F98").Value = "Evan"
$ws.Range("G98").Value = "ES"
$ws.Range("C99").Value = "No"
$ws.Range("D99").Value = "N/A"
$ws.Range("F99").Value = "Evan"
$ws.Range("G99").Value = "ES"
$ws.Range("C100").Value = "No"
$ws.Range("D100").Value = "N/A"
$ws.Range("F100").Value = "Evan"
$ws.Range("G100").Value = "ES"
$ws.Range("C101").Value = "No"
$ws.Range("D101").Value = "N/A"
$ws.Range("F101").Value = "Evan"
$ws.Range("G101").Value = "ES"
$ws.Range("C102").Value = "No"
$ws.Range("D102").Value = "N/A"
$ws.Range("F102").Value = "Evan"
$ws.Range("G102").Value = "ES"
$ws.Range("C103").Value = "No"
$ws.Range("D103").Value = "N/A"
$ws.Range("F103").Value = "Evan"
$ws.Range("G103").Value = "ES"
$ws.Range("C104").Value = "No"
$ws.Range("D104").Value = "N/A"
$ws.Range("F104").Value = "Evan"
$ws.Range("G104").Value = "ES"
$ws.Range("C105").Value = "No"
$ws.Range("D105").Value = "N/A"
$ws.Range("F105").Value = "Evan"
$ws.Range("G105").Value = "ES"
$ws.Range("C106").Value = "No"
$ws.Range("D106").Value = "N/A"
$ws.Range("F106").Value = "Evan"
$ws.Range("G106").Value = "ES"
$ws.Range("C107").Value = "No"
$ws.Range("D107").Value = "N/A"
$ws.Range("F107").Value = "Evan"
$ws.Range("G107").Value = "ES"
$ws.Range("C108").Value = "No"
$ws.Range("D108").Value = "N/A"
$ws.Range("F108").Value = "Evan"
$ws.Range("G108").Value = "ES"
$ws.Range("C109").Value = "No"
$ws.Range("D109").Value = "N/A"
$ws.Range("F109").Value = "Evan"
$ws.Range("G109").Value = "ES"
$ws.Range("C110").Value = "No"
$ws.Range("D110").Value = "N/A"
$ws.Range("F110").Value = "Evan"
$ws.Range("G110").Value = "ES"

# Row heights for wrapped, multi-line text entries
$ws.Rows.Item(92).RowHeight = 30
$ws.Rows.Item(94).RowHeight = 30
$ws.Rows.Item(97).RowHeight = 30
$ws.Rows.Item(99).RowHeight = 45
$ws.Rows.Item(100).RowHeight = 30
$ws.Rows.Item(101).RowHeight = 30
$ws.Rows.Item(102).RowHeight = 30
$ws.Rows.Item(103).RowHeight = 30
$ws.Rows.Item(104).RowHeight = 30
$ws.Rows.Item(105).RowHeight = 45
$ws.Rows.Item(106).RowHeight = 30
$ws.Rows.Item(108).RowHeight = 30
$ws.Rows.Item(110).RowHeight = 30

# --- Extend the blank formatted table by 8 rows (831-838), mirroring row 830 ---
$ws.Range("A830:G830").Copy()
$ws.Range("A831:G838").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A838").Clear()

# --- Update the active view/selection to reflect where data entry left off ---
$ws.Range("A111").Select()
